$wb = $excel.ActiveWorkbook

# ---- Sheet "Metadata" ----
$meta = $wb.Worksheets.Item("Metadata")

# Version
$meta.Range("B3").Value = "2.0.0-sd-202406-matchbox-patch"

# Date
$meta.Range("B8").Value = "2024-06-19T17:47:42+02:00"

# Contact
$meta.Range("B10").Value = "HL7 International - Structured Documents (http://www.hl7.org/Special/committees/structure, structdog@lists.HL7.org)"

# ---- Sheet "Elements" ----
$elements = $wb.Worksheets.Item("Elements")

# Row 2 (PIVL_TS root element) - Definition column (M)
$elements.Range("M2").Value = "A quantity specifying a point on the axis of natural time. A point in time is most often represented as a calendar expression."

# Row 5 (PIVL_TS.operator) - Binding Value Set column (Z)
$elements.Range("Z5").Value = "http://hl7.org/cda/stds/core/ValueSet/CDASetOperator"
